$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.105.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.653.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -0.50%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "218.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5295"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  +1.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.003"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  -0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2609"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06326"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "20.40"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  -3.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07748"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  +0.47%  "

$ws.Range("B12").Value2 = "WrappedEther"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.678.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  +0.81%  "

$ws.Range("B13").Value2 = "Polkadot"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.493"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +1.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.5464"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.0#8130"
$ws.Range("D15").Replace("#", [char]0x2085) | Out-Null
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -1.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "65.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  +0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "26.136.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -0.50%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.537"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -2.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "193.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "10.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.995"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  -1.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.004"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "140.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  +1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.1239"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "7.266"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "16.15"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  +0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.433"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +2.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.05939"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.276"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  -0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.508"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  -3.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.232"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -2.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.545"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -5.27%  "

$ws.Range("B34").Value2 = "HuobiToken"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "2.414"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  -0.04%  "

$ws.Range("B35").Value2 = "ARBITRUM"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.9456"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  -3.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.761"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.5625"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.01609"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  +1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "5.848"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -1.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.8460"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  -1.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.003"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "101.02"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  +1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "1.009.45"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.802.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "56.86"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  -0.58%  "

$ws.Range("B46").Value2 = "Frax"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "1.005"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  +0.15%  "

$ws.Range("B47").Value2 = "BabyDogeCoin"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0#103"
$ws.Range("D47").Replace("#", [char]0x2088) | Out-Null
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  -8.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.4288"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  +1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.470"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.05150"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  -0.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "7.725"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -4.55%  "
